$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Description" column (D) values per corrected buffer distances.
# Order chosen to reproduce the shared-string table order of the target file.
$ws.Range("D5").Value = "250m Buffer"
$ws.Range("D6").Value = "1000m Buffer"
$ws.Range("D4").Value = "500m Buffer (Solar); 1000m Buffer (Wind)"
$ws.Range("D2").Value = "1000m Buffer (Solar); 5000m Buffer (Wind)"
$ws.Range("D3").Value = "1000m Buffer (Solar); 5000m Buffer (Wind)"
$ws.Range("D7").Value = "250m Buffer"
$ws.Range("D8").Value = "250m Buffer"
$ws.Range("D9").Value = "1000m Buffer (Solar); 5000m Buffer (Wind)"

# Rows with shorter, non-wrapping description text shrink back to a single
# line of wrapped text, so their row height auto-fits smaller.
$ws.Rows.Item(5).RowHeight = 29
$ws.Rows.Item(8).RowHeight = 29

# Update the active selection to match the saved workbook state.
$ws.Range("D2").Select()
